$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 424 (shifts existing rows 424:448 down to 426:450),
# inheriting formatting (e.g. the date number format on column D) from the
# row above, same as native Excel "Insert Copied Cells"/"Insert Row" behavior.
$ws.Rows("424:425").Insert()

# New row 424 - "Primera" quality entry for Limache, week of 2021-11-16
$ws.Cells.Item(424, 1).Value = 3
$ws.Cells.Item(424, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(424, 3).Value = "Coquimbo"
$ws.Cells.Item(424, 4).Value = 44516
$ws.Cells.Item(424, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(424, 5).Value = 5
$ws.Cells.Item(424, 6).Value = 100114014
$ws.Cells.Item(424, 7).Value = "Betarraga"
$ws.Cells.Item(424, 8).Value = "Sin especificar"
$ws.Cells.Item(424, 9).Value = "Primera"
$ws.Cells.Item(424, 10).Value = 3200
$ws.Cells.Item(424, 11).Value = 550
$ws.Cells.Item(424, 12).Value = 600
$ws.Cells.Item(424, 13).Value = 575
$ws.Cells.Item(424, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(424, 15).Value = "Limache"
$ws.Cells.Item(424, 16).Value = 144
$ws.Cells.Item(424, 17).Value = 4
$ws.Cells.Item(424, 18).Value = "Hortaliza"

# New row 425 - "Segunda" quality entry for Limache, week of 2021-11-16
$ws.Cells.Item(425, 1).Value = 3
$ws.Cells.Item(425, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(425, 3).Value = "Coquimbo"
$ws.Cells.Item(425, 4).Value = 44516
$ws.Cells.Item(425, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(425, 5).Value = 5
$ws.Cells.Item(425, 6).Value = 100114014
$ws.Cells.Item(425, 7).Value = "Betarraga"
$ws.Cells.Item(425, 8).Value = "Sin especificar"
$ws.Cells.Item(425, 9).Value = "Segunda"
$ws.Cells.Item(425, 10).Value = 1500
$ws.Cells.Item(425, 11).Value = 400
$ws.Cells.Item(425, 12).Value = 400
$ws.Cells.Item(425, 13).Value = 400
$ws.Cells.Item(425, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(425, 15).Value = "Limache"
$ws.Cells.Item(425, 16).Value = 100
$ws.Cells.Item(425, 17).Value = 4
$ws.Cells.Item(425, 18).Value = "Hortaliza"
